# Marksheet regeneration: the grading/marking logic now handles float-valued
# inputs correctly, which changes the computed Right/Wrong/NotAttempt/Marking/
# Total figures and re-lays the "Student Ans / Correct Ans" answer grid from a
# 3-pair-per-row layout (A/B, D/E, G/H) down to a leaner 2-pair layout,
# highlighting matched ("correct") answers in the `correctStyle` (green) face.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ------------------------------------------
# "No." row: label gets the section-title style; Right/Wrong/NotAttempt/Max
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 22
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# "Marking" row: label gets the section-title style; per-question +/- marks
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# "Total" row: label gets the section-title style; total score + fraction
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88/112"

# ---- Answer grid (rows 15-40) ---------------------------------------------
# Drop the third Student-Ans/Correct-Ans column pair (G/H) entirely.
$ws.Range("G15:H40").Clear()

# Drop the second column pair (D/E) for every row past the first block.
$ws.Range("D19:E40").Clear()

# Mirror each "Correct Ans" into the paired "Student Ans" cell wherever the
# student actually answered, formatted with the green `correctStyle` face
# (blank "Student Ans" cells / `normalStyle` formatting are left as-is,
# marking those questions not attempted).
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"
$ws.Range("A24").Value = "Option A"
$ws.Range("A24").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A35").Value = "Option D"
$ws.Range("A35").Style = "correctStyle"
$ws.Range("A37").Value = "Option A"
$ws.Range("A37").Style = "correctStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"
